# Sort demographic sub-tables for comparability:
#  - SEX / RACE / PAY1 sections reordered (SEX now comes before RACE,
#    RACE before PAY1) and every section's rows sorted alphabetically
#  - HOSP_LOCTEACH / HOSP_REGION rows sorted alphabetically
#  - a new INCOME_QRTL section is inserted after HOSP_REGION
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellContent($rowIdx, $colIdx, $text, $bold) {
    $cell = $t.Cell($rowIdx, $colIdx)
    if ($text.Length -eq 0) {
        # Leave/collapse the cell to an empty paragraph (no run)
        $cell.Range.Delete()
    } else {
        $cell.Range.Text = $text
        if ($bold) {
            $full = $cell.Range
            $r = $d.Range($full.Start, $full.Start + $text.Length)
            $r.Font.Bold = $true
        }
    }
}

# Row 1 (title row) is untouched. Drop every other existing row then
# rebuild them in the new, sorted order so formatting stays clean.
for ($i = $t.Rows.Count; $i -ge 2; $i--) {
    $t.Rows.Item($i).Delete()
}

for ($i = 1; $i -le 37; $i++) {
    $t.Rows.Add() | Out-Null
}

Set-CellContent 2 1 'Age' $true
Set-CellContent 2 2 '36.45 +/- 0.06' $false
Set-CellContent 3 1 '>65' $false
Set-CellContent 3 2 '484 (1.21)' $false
Set-CellContent 4 1 'APDRG' $true
Set-CellContent 4 2 '' $false
Set-CellContent 5 1 'APRDRG_Severity > 2' $false
Set-CellContent 5 2 '1,883 (4.70)' $false
Set-CellContent 6 1 'APRDRG_Risk_Mortality > 2' $false
Set-CellContent 6 2 '664 (1.66)' $false
Set-CellContent 7 1 'SEX' $true
Set-CellContent 7 2 '' $false
Set-CellContent 8 1 'Female' $false
Set-CellContent 8 2 '15,475 (38.66)' $false
Set-CellContent 9 1 'Male' $false
Set-CellContent 9 2 '24,385 (60.91)' $false
Set-CellContent 10 1 'Unknown' $false
Set-CellContent 10 2 '173 (0.43)' $false
Set-CellContent 11 1 'RACE' $true
Set-CellContent 11 2 '' $false
Set-CellContent 12 1 'Asian or Pacific Islander' $false
Set-CellContent 12 2 '1,020 (2.55)' $false
Set-CellContent 13 1 'Black' $false
Set-CellContent 13 2 '4,194 (10.48)' $false
Set-CellContent 14 1 'Hispanic' $false
Set-CellContent 14 2 '13,494 (33.71)' $false
Set-CellContent 15 1 'Native American' $false
Set-CellContent 15 2 '213 (0.53)' $false
Set-CellContent 16 1 'Other' $false
Set-CellContent 16 2 '2,157 (5.39)' $false
Set-CellContent 17 1 'Unknown' $false
Set-CellContent 17 2 '1,911 (4.77)' $false
Set-CellContent 18 1 'White' $false
Set-CellContent 18 2 '17,044 (42.57)' $false
Set-CellContent 19 1 'PAY1' $true
Set-CellContent 19 2 '' $false
Set-CellContent 20 1 'Self-pay' $false
Set-CellContent 20 2 '40,033 (100.00)' $false
Set-CellContent 21 1 'HOSP_LOCTEACH' $true
Set-CellContent 21 2 '' $false
Set-CellContent 22 1 'Rural' $false
Set-CellContent 22 2 '4,698 (11.74)' $false
Set-CellContent 23 1 'Urban nonteaching' $false
Set-CellContent 23 2 '15,200 (37.97)' $false
Set-CellContent 24 1 'Urban teaching' $false
Set-CellContent 24 2 '20,135 (50.30)' $false
Set-CellContent 25 1 'HOSP_REGION' $true
Set-CellContent 25 2 '' $false
Set-CellContent 26 1 'Midwest' $false
Set-CellContent 26 2 '5,816 (14.53)' $false
Set-CellContent 27 1 'Northeast' $false
Set-CellContent 27 2 '5,398 (13.48)' $false
Set-CellContent 28 1 'South' $false
Set-CellContent 28 2 '21,211 (52.98)' $false
Set-CellContent 29 1 'West' $false
Set-CellContent 29 2 '7,608 (19.00)' $false
Set-CellContent 30 1 'INCOME_QRTL' $true
Set-CellContent 30 2 '' $false
Set-CellContent 31 1 '1' $false
Set-CellContent 31 2 '14,830 (37.04)' $false
Set-CellContent 32 1 '2' $false
Set-CellContent 32 2 '10,937 (27.32)' $false
Set-CellContent 33 1 '3' $false
Set-CellContent 33 2 '8,942 (22.34)' $false
Set-CellContent 34 1 '4' $false
Set-CellContent 34 2 '5,324 (13.30)' $false
Set-CellContent 35 1 'SSI' $false
Set-CellContent 35 2 '178 (0.44)' $false
Set-CellContent 36 1 'PROLONGED_LOS' $false
Set-CellContent 36 2 '6,002 (14.99)' $false
Set-CellContent 37 1 'DIED' $false
Set-CellContent 37 2 '25 (0.06)' $false
Set-CellContent 38 1 'OR_RETURN' $false
Set-CellContent 38 2 '9,227 (23.05)' $false

Write-Output ("Final row count=" + $t.Rows.Count)
